$d = $word.ActiveDocument

$replacements = @(
    @{old = "648×8="; new = "723×7="},
    @{old = "625×9="; new = "731×5="},
    @{old = "365×2="; new = "504×6="},
    @{old = "199×5="; new = "656×3="},
    @{old = "384×5="; new = "557×4="},
    @{old = "919×7="; new = "661×9="},
    @{old = "558×8="; new = "771×8="},
    @{old = "286×5="; new = "294×7="},
    @{old = "715×3="; new = "331×7="},
    @{old = "572×2="; new = "143×2="},
    @{old = "304×8="; new = "819×2="},
    @{old = "952×2="; new = "661×7="},
    @{old = "806×7="; new = "133×8="},
    @{old = "733×7="; new = "924×2="},
    @{old = "216×5="; new = "692×5="},
    @{old = "307×7="; new = "400×9="},
    @{old = "976×7="; new = "582×7="},
    @{old = "294×2="; new = "681×6="},
    @{old = "228×6="; new = "500×4="},
    @{old = "991×8="; new = "312×3="},
    @{old = "812×2="; new = "839×7="},
    @{old = "762×6="; new = "831×5="},
    @{old = "267×7="; new = "575×3="},
    @{old = "445×3="; new = "700×3="},
    @{old = "777×5="; new = "745×7="}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $pair.new, 2)
}
